$wb = $excel.ActiveWorkbook

# Sheet "展览": update 想去人数 (interest count) values for two events
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F12").Value = 11976
$wsExhibit.Range("F13").Value = 5436

# Sheet "全部类型": same two events also listed here, update accordingly
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F14").Value = 11976
$wsAll.Range("F16").Value = 5436
